# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 2-4, column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1002
$wsExhibit.Range("F3").Value = 2061
$wsExhibit.Range("F4").Value = 454

# Sheet "全部类型" - rows 4-6, column F (same events repeated further down the list)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1002
$wsAll.Range("F5").Value = 2061
$wsAll.Range("F6").Value = 454
